# Update crypto price (D) and 1h volume % (E) columns with refreshed quotes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Price = "319.29"; Volume = "3.68%" }
    @{ Row = 3; Price = "48.93"; Volume = "11.83%" }
    @{ Row = 4; Price = "5.317"; Volume = "4.11%" }
    @{ Row = 5; Price = "0.07961"; Volume = "-0.16%" }
    @{ Row = 6; Price = "4.605"; Volume = "3.48%" }
    @{ Row = 7; Price = "1.326"; Volume = "23.61%" }
    @{ Row = 8; Price = "1.643"; Volume = "1.62%" }
    @{ Row = 9; Price = "0.1247"; Volume = "-3.08%" }
    @{ Row = 10; Price = "0.1955"; Volume = "3.51%" }
    @{ Row = 11; Price = "0.09457"; Volume = "3.22%" }
    @{ Row = 12; Price = "0.04581"; Volume = "9.76%" }
    @{ Row = 13; Price = "0.1048"; Volume = "0.68%" }
    @{ Row = 14; Price = "0.001306"; Volume = "1.68%" }
    @{ Row = 15; Price = "0.04217"; Volume = "1.63%" }
    @{ Row = 16; Price = "0.006010"; Volume = "5.37%" }
    @{ Row = 17; Price = "3.337"; Volume = "-0.09%" }
    @{ Row = 18; Price = "2.439"; Volume = "2.59%" }
    @{ Row = 19; Price = "0.3466"; Volume = "3.41%" }
    @{ Row = 20; Price = "8.103"; Volume = "1.59%" }
    @{ Row = 21; Price = "0.1410"; Volume = "2.61%" }
    @{ Row = 22; Price = "0.3073"; Volume = "10.17%" }
    @{ Row = 23; Price = "0.001301"; Volume = "3.19%" }
    @{ Row = 24; Price = "0.004182"; Volume = "-7.56%" }
    @{ Row = 25; Price = "0.0001355"; Volume = "2.45%" }
    @{ Row = 26; Price = "0.0003569"; Volume = "-95.18%" }
    @{ Row = 38; Price = "0.02627"; Volume = "-1.95%" }
    @{ Row = 39; Price = "0.05852"; Volume = "9.06%" }
    @{ Row = 40; Price = "0.01038"; Volume = "86.82%" }
    @{ Row = 41; Price = "0.008030"; Volume = "3.41%" }
    @{ Row = 42; Price = "0.1453"; Volume = "3.52%" }
    @{ Row = 43; Price = "0.007578"; Volume = "4.20%" }
    @{ Row = 44; Price = "0.007947"; Volume = "-4.55%" }
    @{ Row = 45; Price = "0.3201"; Volume = "4.55%" }
    @{ Row = 46; Price = "0.00007063"; Volume = "6.31%" }
    @{ Row = 47; Price = "0.00000000756"; Volume = "2.81%" }
    @{ Row = 48; Price = "0.05598"; Volume = "-9.05%" }
    @{ Row = 49; Price = "0.004034"; Volume = "2.12%" }
    @{ Row = 50; Price = "0.00002118"; Volume = "2.81%" }
    @{ Row = 51; Price = "0.0002017"; Volume = "2.81%" }
)

foreach ($u in $updates) {
    $dCell = $ws.Range("D" + $u.Row)
    $dCell.Value = "'" + $u.Price
    $dCell.Style = "Normal"

    $eCell = $ws.Range("E" + $u.Row)
    $eCell.Value = "'" + $u.Volume
    $eCell.Style = "Normal"
}

